$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = '[-0.1534658717471018, 13.558235014520294]'
$ws.Range("N2").Value = [double]"0.05512751251285741"
$ws.Range("O2").Value = [double]"0.05512751251285741"
$ws.Range("Q2").Value = '[-3.018947895341543, -0.22642109215061534]'
$ws.Range("R2").Value = [double]"0.02373835271602731"
$ws.Range("S2").Value = [double]"0.02373835271602731"
$ws.Range("U2").Value = '[5.145471555984354, 13.104807805090612]'
$ws.Range("V2").Value = [double]"3.240845782781498e-05"
$ws.Range("W2").Value = [double]"3.240845782781498e-05"
$ws.Range("Y2").Value = [double]"0.9192792792792979"
$ws.Range("Z2").Value = [double]"12.25705705705732"
$ws.Range("M3").Value = '[-1.171715976062032, 14.22280365671034]'
$ws.Range("N3").Value = [double]"0.09461958417285343"
$ws.Range("O3").Value = [double]"0.09461958417285343"
$ws.Range("Q3").Value = '[-3.660474323101621, -0.08805264694746207]'
$ws.Range("R3").Value = [double]"0.04014509924346088"
$ws.Range("S3").Value = [double]"0.04014509924346088"
$ws.Range("U3").Value = '[5.762308671957012, 14.431713615624144]'
$ws.Range("V3").Value = [double]"2.549751884362017e-05"
$ws.Range("W3").Value = [double]"2.549751884362017e-05"
$ws.Range("Y3").Value = [double]"0.3574974974975049"
$ws.Range("Z3").Value = [double]"14.861681681682"
$ws.Range("M4").Value = '[-1.035746136892545, 14.104000095998062]'
$ws.Range("N4").Value = [double]"0.08895706566849415"
$ws.Range("O4").Value = [double]"0.08895706566849415"
$ws.Range("Q4").Value = '[-4.1762112552224675, -0.07547369738353815]'
$ws.Range("R4").Value = [double]"0.04246526634611536"
$ws.Range("S4").Value = [double]"0.04246526634611536"
$ws.Range("U4").Value = '[4.932300433674429, 12.851858251582629]'
$ws.Range("V4").Value = [double]"4.419879816786043e-05"
$ws.Range("W4").Value = [double]"4.419879816786043e-05"
$ws.Range("Y4").Value = [double]"0.306426426426432"
$ws.Range("Z4").Value = [double]"16.95559559559596"
$ws.Range("M5").Value = '[-0.5113816428352216, 13.912683726339312]'
$ws.Range("N5").Value = [double]"0.06781479712738681"
$ws.Range("O5").Value = [double]"0.06781479712738681"
$ws.Range("Q5").Value = '[-4.603895540395853, -0.754736973835386]'
$ws.Range("R5").Value = [double]"0.007425163044326943"
$ws.Range("S5").Value = [double]"0.007425163044326943"
$ws.Range("U5").Value = '[4.836214940865591, 12.760749160364103]'
$ws.Range("V5").Value = [double]"5.203883876658821e-05"
$ws.Range("W5").Value = [double]"5.203883876658821e-05"
$ws.Range("Y5").Value = [double]"3.064264264264332"
$ws.Range("Z5").Value = [double]"18.69201201201242"
$ws.Range("M6").Value = '[-1.1795463000097346, 14.140296447739438]'
$ws.Range("N6").Value = [double]"0.09528694714022223"
$ws.Range("O6").Value = [double]"0.09528694714022223"
$ws.Range("Q6").Value = '[-4.742263985599006, -1.0440528138056178]'
$ws.Range("R6").Value = [double]"0.002889372130075563"
$ws.Range("S6").Value = [double]"0.002889372130075563"
$ws.Range("U6").Value = '[5.501218286381442, 14.229170393012975]'
$ws.Range("V6").Value = [double]"4.006896393615733e-05"
$ws.Range("W6").Value = [double]"4.006896393615733e-05"
$ws.Range("Y6").Value = [double]"4.238898898898996"
$ws.Range("Z6").Value = [double]"19.25379379379421"
$ws.Range("M7").Value = '[-0.08929624896798671, 13.54765501085797]'
$ws.Range("N7").Value = [double]"0.05294728186435527"
$ws.Range("O7").Value = [double]"0.05294728186435527"
$ws.Range("U7").Value = '[4.880393862041759, 12.885509060865067]'
$ws.Range("V7").Value = [double]"5.246372511669151e-05"
$ws.Range("W7").Value = [double]"5.246372511669151e-05"
$ws.Range("M8").Value = '[-1.1008571175800785, 14.179897686262027]'
$ws.Range("N8").Value = [double]"0.09159010742149509"
$ws.Range("O8").Value = [double]"0.09159010742149509"
$ws.Range("Q8").Value = '[0.05660527303765317, 5.163658795990431]'
$ws.Range("R8").Value = [double]"0.04533471763134544"
$ws.Range("S8").Value = [double]"0.04533471763134544"
$ws.Range("U8").Value = '[4.942190781597072, 12.89334933519241]'
$ws.Range("V8").Value = [double]"4.491566888709109e-05"
$ws.Range("W8").Value = [double]"4.491566888709109e-05"
$ws.Range("Y8").Value = [double]"4.545325325325422"
$ws.Range("Z8").Value = [double]"25.28018018018073"
$ws.Range("M9").Value = '[-1.841656250441396, 14.788252056942277]'
$ws.Range("N9").Value = [double]"0.1238857682086514"
$ws.Range("O9").Value = [double]"0.1238857682086514"
$ws.Range("Q9").Value = '[-0.792473822527155, 5.484422009870468]'
$ws.Range("R9").Value = [double]"0.1391749614514639"
$ws.Range("S9").Value = [double]"0.1391749614514639"
$ws.Range("U9").Value = '[5.725174616208108, 14.328957699983008]'
$ws.Range("V9").Value = [double]"2.524505588397119e-05"
$ws.Range("W9").Value = [double]"2.524505588397119e-05"
$ws.Range("Y9").Value = [double]"3.243013013013087"
$ws.Range("Z9").Value = [double]"28.7274774774781"
$ws.Range("M10").Value = '[-1.1759274130162325, 14.808223773819012]'
$ws.Range("N10").Value = [double]"0.09271647090797575"
$ws.Range("O10").Value = [double]"0.09271647090797575"
$ws.Range("Q10").Value = '[0.16981581911296217, 3.767395394394967]'
$ws.Range("R10").Value = [double]"0.03266342707563363"
$ws.Range("S10").Value = [double]"0.03266342707563363"
$ws.Range("U10").Value = '[5.552011429438236, 13.942908931623466]'
$ws.Range("V10").Value = [double]"2.652786516055272e-05"
$ws.Range("W10").Value = [double]"2.652786516055272e-05"
$ws.Range("Y10").Value = [double]"9.849849849850013"
$ws.Range("Z10").Value = [double]"23.93513513513553"
$ws.Range("M11").Value = '[-0.7870469925739307, 14.423962249665685]'
$ws.Range("N11").Value = [double]"0.07766091994408875"
$ws.Range("O11").Value = [double]"0.07766091994408875"
$ws.Range("Q11").Value = '[0.20755266780473214, 3.528395352680427]'
$ws.Range("R11").Value = [double]"0.02831859753591637"
$ws.Range("S11").Value = [double]"0.02831859753591637"
$ws.Range("U11").Value = '[5.531141472571845, 13.913708022395573]'
$ws.Range("V11").Value = [double]"2.71768716975096e-05"
$ws.Range("W11").Value = [double]"2.71768716975096e-05"
$ws.Range("Y11").Value = [double]"10.78558558558577"
$ws.Range("Z11").Value = [double]"23.78738738738777"
$ws.Range("M12").Value = '[0.8249184445311037, 12.816618551531068]'
$ws.Range("N12").Value = [double]"0.02668654386422853"
$ws.Range("O12").Value = [double]"0.02668654386422853"
$ws.Range("Q12").Value = '[0.5471843060306547, 3.012658420559582]'
$ws.Range("R12").Value = [double]"0.005628918868783694"
$ws.Range("S12").Value = [double]"0.005628918868783694"
$ws.Range("U12").Value = '[4.303095213695277, 11.254900540570098]'
$ws.Range("V12").Value = [double]"4.645425796745428e-05"
$ws.Range("W12").Value = [double]"4.645425796745428e-05"
$ws.Range("Y12").Value = [double]"12.80480480480501"
$ws.Range("Z12").Value = [double]"22.45765765765803"
$ws.Range("M13").Value = '[0.3157668208163269, 12.969140937830872]'
$ws.Range("N13").Value = [double]"0.04003383262178528"
$ws.Range("O13").Value = [double]"0.04003383262178528"
$ws.Range("Q13").Value = '[0.40881586082750054, 2.8742899753564277]'
$ws.Range("R13").Value = [double]"0.01019575836543307"
$ws.Range("S13").Value = [double]"0.01019575836543307"
$ws.Range("U13").Value = '[4.328821352819722, 11.610562081310611]'
$ws.Range("V13").Value = [double]"6.390647196985455e-05"
$ws.Range("W13").Value = [double]"6.390647196985455e-05"
$ws.Range("Y13").Value = [double]"13.34654654654677"
$ws.Range("Z13").Value = [double]"22.99939939939978"
$ws.Range("M14").Value = '[0.08414062820718904, 13.007309660466285]'
$ws.Range("N14").Value = [double]"0.04721150351705639"
$ws.Range("O14").Value = [double]"0.04721150351705639"
$ws.Range("Q14").Value = '[0.14465791998511524, 2.8617110257925047]'
$ws.Range("R14").Value = [double]"0.03088249893073725"
$ws.Range("S14").Value = [double]"0.03088249893073725"
$ws.Range("U14").Value = '[4.375790744684283, 11.815253054773876]'
$ws.Range("V14").Value = [double]"6.933169199685452e-05"
$ws.Range("W14").Value = [double]"6.933169199685452e-05"
$ws.Range("Y14").Value = [double]"13.39579579579601"
$ws.Range("Z14").Value = [double]"24.03363363363403"
$ws.Range("M15").Value = '[-0.34270535968490456, 13.473684344901992]'
$ws.Range("N15").Value = [double]"0.06196835135467404"
$ws.Range("O15").Value = [double]"0.06196835135467404"
$ws.Range("Q15").Value = '[-0.09434212172942225, 2.8491320762285817]'
$ws.Range("R15").Value = [double]"0.06589800617748121"
$ws.Range("S15").Value = [double]"0.06589800617748121"
$ws.Range("U15").Value = '[4.4637983409665045, 12.101062981341869]'
$ws.Range("V15").Value = [double]"7.274049283467754e-05"
$ws.Range("W15").Value = [double]"7.274049283467754e-05"
$ws.Range("Y15").Value = [double]"13.44504504504527"
$ws.Range("Z15").Value = [double]"24.96936936936978"
